$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 2131.6667
$ws.Range("J121").Value = 2775
$ws.Range("L121").Value = 8325
$ws.Range("N121").Value = -11819
$ws.Range("H137").Value = 1694.6279
$ws.Range("I137").Value = 1539.6666
$ws.Range("J137").Value = 2052.2307
$ws.Range("K137").Value = 4618.9998
$ws.Range("L137").Value = 6156.6921
$ws.Range("M137").Value = -2068.9998
$ws.Range("N137").Value = -11256.6921
$ws.Range("H138").Value = 5822.6055
$ws.Range("I138").Value = 1412
$ws.Range("J138").Value = 7858.269
$ws.Range("K138").Value = 4236
$ws.Range("L138").Value = 23574.807
$ws.Range("M138").Value = 904
$ws.Range("N138").Value = -33854.807

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8094.409
$ws.Range("I32").Value = 3806.1296
$ws.Range("K32").Value = 3806.1296
$ws.Range("M32").Value = -3519.1296
$ws.Range("H41").Value = 3518.6667
$ws.Range("I41").Value = 3518.6667
$ws.Range("K41").Value = 3518.6667
$ws.Range("M41").Value = -3104.6667
$ws.Range("H61").Value = 8427.277
$ws.Range("I61").Value = 8814.066000000001
$ws.Range("J61").Value = 6493.3335
$ws.Range("K61").Value = 8814.066000000001
$ws.Range("L61").Value = 6493.3335
$ws.Range("M61").Value = -8602.066000000001
$ws.Range("N61").Value = -6917.3335
$ws.Range("H74").Value = 1739.2609
$ws.Range("I74").Value = 1467.4706
$ws.Range("J74").Value = 2509.3333
$ws.Range("K74").Value = 1467.4706
$ws.Range("L74").Value = 2509.3333
$ws.Range("M74").Value = -593.4706000000001
$ws.Range("N74").Value = -4257.3333
$ws.Range("H77").Value = 1739.2609
$ws.Range("I77").Value = 1467.4706
$ws.Range("J77").Value = 2509.3333
$ws.Range("K77").Value = 7337.353000000001
$ws.Range("L77").Value = 12546.6665
$ws.Range("M77").Value = -2969.353000000001
$ws.Range("N77").Value = -21282.6665
$ws.Range("H97").Value = 1642
$ws.Range("I97").Value = 1646.3636
$ws.Range("J97").Value = 1630
$ws.Range("K97").Value = 1646.3636
$ws.Range("L97").Value = 1630
$ws.Range("M97").Value = -1150.3636
$ws.Range("N97").Value = -2622
$ws.Range("H122").Value = 3664331.8
$ws.Range("I122").Value = 3664331.8
$ws.Range("K122").Value = 10992995.4
$ws.Range("M122").Value = -10990545.4
$ws.Range("H124").Value = 25214.5
$ws.Range("J124").Value = 25214.5
$ws.Range("L124").Value = 25214.5
$ws.Range("N124").Value = -35034.5
$ws.Range("H132").Value = 4415
$ws.Range("I132").Value = 2209.4167
$ws.Range("J132").Value = 7723.375
$ws.Range("K132").Value = 6628.250100000001
$ws.Range("L132").Value = 23170.125
$ws.Range("M132").Value = -4098.250100000001
$ws.Range("N132").Value = -28230.125
$ws.Range("H136").Value = 8427.277
$ws.Range("I136").Value = 8814.066000000001
$ws.Range("J136").Value = 6493.3335
$ws.Range("K136").Value = 26442.198
$ws.Range("L136").Value = 19480.0005
$ws.Range("M136").Value = -23892.198
$ws.Range("N136").Value = -24580.0005

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("N9").Value = $null
$ws.Range("H44").Value = 10000
$ws.Range("J44").Value = 10000
$ws.Range("L44").Value = 10000
$ws.Range("N44").Value = -10994
$ws.Range("H49").Value = 11900
$ws.Range("J49").Value = 11900
$ws.Range("L49").Value = 11900
$ws.Range("N49").Value = -12378
$ws.Range("H107").Value = 1063.625
$ws.Range("I107").Value = 951.8333
$ws.Range("J107").Value = 1399
$ws.Range("K107").Value = 951.8333
$ws.Range("L107").Value = 1399
$ws.Range("M107").Value = 968.1667
$ws.Range("N107").Value = -5239
$ws.Range("H134").Value = 13518335
$ws.Range("I134").Value = 18524022
$ws.Range("J134").Value = 2980
$ws.Range("K134").Value = 55572066
$ws.Range("L134").Value = 8940
$ws.Range("M134").Value = -55569531
$ws.Range("N134").Value = -14010

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4014.7937
$ws.Range("I31").Value = 1575.3334
$ws.Range("K31").Value = 1575.3334
$ws.Range("M31").Value = -1280.3334
$ws.Range("H34").Value = 4014.7937
$ws.Range("I34").Value = 1575.3334
$ws.Range("K34").Value = 1575.3334
$ws.Range("M34").Value = -1373.3334

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3217.5334
$ws.Range("I80").Value = 2994.3333
$ws.Range("K80").Value = 2994.3333
$ws.Range("M80").Value = -1996.3333
$ws.Range("H83").Value = 3217.5334
$ws.Range("I83").Value = 2994.3333
$ws.Range("K83").Value = 14971.6665
$ws.Range("M83").Value = -9979.666499999999
$ws.Range("H111").Value = 20000
$ws.Range("J111").Value = 20000
$ws.Range("L111").Value = 20000
$ws.Range("N111").Value = -26134
$ws.Range("H126").Value = 5203.9033
$ws.Range("I126").Value = 7001.3687
$ws.Range("J126").Value = 2357.9167
$ws.Range("K126").Value = 21004.1061
$ws.Range("L126").Value = 7073.750100000001
$ws.Range("M126").Value = -18534.1061
$ws.Range("N126").Value = -12013.7501
$ws.Range("H132").Value = 5803.5835
$ws.Range("I132").Value = 10703.75
$ws.Range("J132").Value = 3353.5
$ws.Range("K132").Value = 32111.25
$ws.Range("L132").Value = 10060.5
$ws.Range("M132").Value = -29581.25
$ws.Range("N132").Value = -15120.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1414
$ws.Range("I100").Value = 1002.4
$ws.Range("J100").Value = 2100
$ws.Range("K100").Value = 1002.4
$ws.Range("L100").Value = 2100
$ws.Range("M100").Value = -461.4
$ws.Range("N100").Value = -3182
$ws.Range("H122").Value = 2396980
$ws.Range("I122").Value = 3761221.8
$ws.Range("J122").Value = 668940.7
$ws.Range("K122").Value = 11283665.4
$ws.Range("L122").Value = 2006822.1
$ws.Range("M122").Value = -11281215.4
$ws.Range("N122").Value = -2011722.1
$ws.Range("H136").Value = 10123.303
$ws.Range("I136").Value = 6829.875
$ws.Range("J136").Value = 18905.777
$ws.Range("K136").Value = 20489.625
$ws.Range("L136").Value = 56717.33099999999
$ws.Range("M136").Value = -17939.625
$ws.Range("N136").Value = -61817.33099999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 998.2432
$ws.Range("I113").Value = 918.63635
$ws.Range("J113").Value = 1115
$ws.Range("K113").Value = 2755.90905
$ws.Range("L113").Value = 3345
$ws.Range("M113").Value = -585.9090500000002
$ws.Range("N113").Value = -7685
